$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Langth", $true, $false, $false, $false, $false, $true, 1, $false, "Length", 2)
